# Apply updated TPM-derived values to the Gdf11-Acvr1b LR-pair worksheet
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 0.3288063333333333
$ws.Range("H2").Value = 0.9864189999999999
$ws.Range("I2").Value = 0.05575527297994041
$ws.Range("J2").Value = 0.05575527297994041
$ws.Range("M2").Value = 1.522526333333333
$ws.Range("N2").Value = 4.567579
$ws.Range("O2").Value = 0.2115373313282365
$ws.Range("P2").Value = 0.2115373313282365
$ws.Range("Q2").Value = 0.5006163010667778
$ws.Range("R2").Value = 4.505546709601
$ws.Range("S2").Value = 0.01179432165365393
$ws.Range("T2").Value = 0.01179432165365393
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 0.6666666666666666
$ws.Range("G3").Value = 0.3288063333333333
$ws.Range("H3").Value = 0.9864189999999999
$ws.Range("I3").Value = 0.05575527297994041
$ws.Range("J3").Value = 0.05575527297994041
$ws.Range("O3").Value = 0.4376697219060474
$ws.Range("P3").Value = 0.4376697219060474
$ws.Range("Q3").Value = 1.035772721031222
$ws.Range("R3").Value = 9.321954489281001
$ws.Range("S3").Value = 0.02440239481992628
$ws.Range("T3").Value = 0.02440239481992628
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 0.6666666666666666
$ws.Range("G4").Value = 0.3288063333333333
$ws.Range("H4").Value = 0.9864189999999999
$ws.Range("I4").Value = 0.05575527297994041
$ws.Range("J4").Value = 0.05575527297994041
$ws.Range("M4").Value = 2.524809666666667
$ws.Range("N4").Value = 7.574429
$ws.Range("O4").Value = 0.3507929467657161
$ws.Range("P4").Value = 0.3507929467657162
$ws.Range("Q4").Value = 0.8301734088612222
$ws.Range("R4").Value = 7.471560679751
$ws.Range("S4").Value = 0.01955855650636021
$ws.Range("T4").Value = 0.01955855650636021
$ws.Range("I5").Value = 0.3115445049245869
$ws.Range("J5").Value = 0.3115445049245869
$ws.Range("M5").Value = 1.522526333333333
$ws.Range("N5").Value = 4.567579
$ws.Range("O5").Value = 0.2115373313282365
$ws.Range("P5").Value = 0.2115373313282365
$ws.Range("Q5").Value = 2.797300584092556
$ws.Range("R5").Value = 25.175705256833
$ws.Range("S5").Value = 0.06590329316172373
$ws.Range("T5").Value = 0.06590329316172373
$ws.Range("I6").Value = 0.3115445049245869
$ws.Range("J6").Value = 0.3115445049245869
$ws.Range("O6").Value = 0.4376697219060474
$ws.Range("P6").Value = 0.4376697219060474
$ws.Range("S6").Value = 0.1363535968317011
$ws.Range("T6").Value = 0.1363535968317011
$ws.Range("I7").Value = 0.3115445049245869
$ws.Range("J7").Value = 0.3115445049245869
$ws.Range("M7").Value = 2.524809666666667
$ws.Range("N7").Value = 7.574429
$ws.Range("O7").Value = 0.3507929467657161
$ws.Range("P7").Value = 0.3507929467657162
$ws.Range("Q7").Value = 4.638771363531444
$ws.Range("R7").Value = 41.748942271783
$ws.Range("S7").Value = 0.109287614931162
$ws.Range("T7").Value = 0.109287614931162
$ws.Range("G8").Value = 3.731231666666667
$ws.Range("H8").Value = 11.193695
$ws.Range("I8").Value = 0.6327002220954728
$ws.Range("J8").Value = 0.6327002220954728
$ws.Range("M8").Value = 1.522526333333333
$ws.Range("N8").Value = 4.567579
$ws.Range("O8").Value = 0.2115373313282365
$ws.Range("P8").Value = 0.2115373313282365
$ws.Range("Q8").Value = 5.680898468267224
$ws.Range("R8").Value = 51.12808621440501
$ws.Range("S8").Value = 0.1338397165128588
$ws.Range("T8").Value = 0.1338397165128588
$ws.Range("G9").Value = 3.731231666666667
$ws.Range("H9").Value = 11.193695
$ws.Range("I9").Value = 0.6327002220954728
$ws.Range("J9").Value = 0.6327002220954728
$ws.Range("O9").Value = 0.4376697219060474
$ws.Range("P9").Value = 0.4376697219060474
$ws.Range("Q9").Value = 11.75375162942278
$ws.Range("R9").Value = 105.783764664805
$ws.Range("S9").Value = 0.27691373025442
$ws.Range("T9").Value = 0.27691373025442
$ws.Range("G10").Value = 3.731231666666667
$ws.Range("H10").Value = 11.193695
$ws.Range("I10").Value = 0.6327002220954728
$ws.Range("J10").Value = 0.6327002220954728
$ws.Range("M10").Value = 2.524809666666667
$ws.Range("N10").Value = 7.574429
$ws.Range("O10").Value = 0.3507929467657161
$ws.Range("P10").Value = 0.3507929467657162
$ws.Range("Q10").Value = 9.420649780572779
$ws.Range("R10").Value = 84.78584802515502
$ws.Range("S10").Value = 0.221946775328194
$ws.Range("T10").Value = 0.221946775328194
